$wb = $excel.ActiveWorkbook

# --- optimization_parameters: drop the obsolete "Sheet" row (old row 16:
#     A="Sheet", B=3, C=4). Deleting the whole row also removes the now
#     -unreferenced "Sheet" shared string and shifts every row below it
#     (and every shared-string / style index after it) up/down by one,
#     matching the diff automatically. ---
$wsOpt = $wb.Worksheets.Item("optimization_parameters")
$wsOpt.Rows.Item(16).Delete()

# After the delete, the row that used to be 17 (simulation_timepoints) is
# now row 16 and is the one Excel left selected as a whole-row selection.
$wsOpt.Rows.Item(16).Select()

# --- network: update the lingering selection left over from editing ---
$wsNetwork = $wb.Worksheets.Item("network")
$wsNetwork.Range("A41").Select()

# --- make degradation_rates the active/visible tab again, keeping its
#     existing B1 selection intact ---
$wsDeg = $wb.Worksheets.Item("degradation_rates")
$wsDeg.Activate()
